$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7170026666666667
$ws.Range("H2").Value = 2.151008
$ws.Range("I2").Value = 0.02953485643833859
$ws.Range("J2").Value = 0.02953485643833859
$ws.Range("M2").Value = 201.098592
$ws.Range("N2").Value = 603.295776
$ws.Range("O2").Value = 0.7918622805845071
$ws.Range("P2").Value = 0.791862280584507
$ws.Range("Q2").Value = 144.188226726912
$ws.Range("R2").Value = 1297.694040542208
$ws.Range("S2").Value = 0.02338753877599881
$ws.Range("T2").Value = 0.02338753877599881
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7170026666666667
$ws.Range("H3").Value = 2.151008
$ws.Range("I3").Value = 0.02953485643833859
$ws.Range("J3").Value = 0.02953485643833859
$ws.Range("O3").Value = 0.1414593902976603
$ws.Range("P3").Value = 0.1414593902976603
$ws.Range("Q3").Value = 25.75798739376
$ws.Range("R3").Value = 231.82188654384
$ws.Range("S3").Value = 0.004177982784296303
$ws.Range("T3").Value = 0.004177982784296303
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7170026666666667
$ws.Range("H4").Value = 2.151008
$ws.Range("I4").Value = 0.02953485643833859
$ws.Range("J4").Value = 0.02953485643833859
$ws.Range("O4").Value = 0.0666783291178327
$ws.Range("P4").Value = 0.06667832911783268
$ws.Range("Q4").Value = 12.14129056572444
$ws.Range("R4").Value = 109.27161509152
$ws.Range("S4").Value = 0.00196933487804348
$ws.Range("T4").Value = 0.00196933487804348
$ws.Range("I5").Value = 0.4970672037825566
$ws.Range("J5").Value = 0.4970672037825566
$ws.Range("M5").Value = 201.098592
$ws.Range("N5").Value = 603.295776
$ws.Range("O5").Value = 0.7918622805845071
$ws.Range("P5").Value = 0.791862280584507
$ws.Range("Q5").Value = 2426.666228330688
$ws.Range("R5").Value = 21839.99605497619
$ws.Range("S5").Value = 0.3936087695910192
$ws.Range("T5").Value = 0.3936087695910191
$ws.Range("I6").Value = 0.4970672037825566
$ws.Range("J6").Value = 0.4970672037825566
$ws.Range("O6").Value = 0.1414593902976603
$ws.Range("P6").Value = 0.1414593902976603
$ws.Range("Q6").Value = 433.5030642729899
$ws.Range("R6").Value = 3901.52757845691
$ws.Range("S6").Value = 0.07031482358404331
$ws.Range("T6").Value = 0.07031482358404331
$ws.Range("I7").Value = 0.4970672037825566
$ws.Range("J7").Value = 0.4970672037825566
$ws.Range("O7").Value = 0.0666783291178327
$ws.Range("P7").Value = 0.06667832911783268
$ws.Range("Q7").Value = 204.3360990907755
$ws.Range("S7").Value = 0.03314361060749412
$ws.Range("T7").Value = 0.03314361060749411
$ws.Range("I8").Value = 0.4733979397791048
$ws.Range("J8").Value = 0.4733979397791048
$ws.Range("M8").Value = 201.098592
$ws.Range("N8").Value = 603.295776
$ws.Range("O8").Value = 0.7918622805845071
$ws.Range("P8").Value = 0.791862280584507
$ws.Range("Q8").Value = 2311.113636710209
$ws.Range("R8").Value = 20800.02273039187
$ws.Range("S8").Value = 0.3748659722174891
$ws.Range("T8").Value = 0.3748659722174891
$ws.Range("I9").Value = 0.4733979397791048
$ws.Range("J9").Value = 0.4733979397791048
$ws.Range("O9").Value = 0.1414593902976603
$ws.Range("P9").Value = 0.1414593902976603
$ws.Range("S9").Value = 0.06696658392932066
$ws.Range("T9").Value = 0.06696658392932066
$ws.Range("I10").Value = 0.4733979397791048
$ws.Range("J10").Value = 0.4733979397791048
$ws.Range("O10").Value = 0.0666783291178327
$ws.Range("P10").Value = 0.06667832911783268
$ws.Range("S10").Value = 0.03156538363229509
$ws.Range("T10").Value = 0.03156538363229509
